$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# --- Add new row 3 (copy formatting from existing row 2) ---
$ws.Range("A2:E2").Copy()
$ws.Range("A3:E3").PasteSpecial(-4122)

# --- Add new row 4 (copy formatting from new row 3) ---
$ws.Range("A3:E3").Copy()
$ws.Range("A4:E4").PasteSpecial(-4122)

# --- Fill in the values (order matters for shared-string allocation) ---
$ws.Range("C3").Value = "Two Pointers"
$ws.Range("A3").Value = "121. Best Time to Buy and Sell Stock"
$ws.Range("D3").Value = "Kadane's is optimal,  but 2 pointers is best for intereview setting. Updating least so far is also intuitive."
$ws.Range("E3").Value = "https://leetcode.com/problems/best-time-to-buy-and-sell-stock/solutions/1735493/java-c-best-ever-explanation-could-possible/ "
$ws.Range("B3").Value = "Easy"

$ws.Range("D4").Value = "Remember the recursive relation argmax(dfs(i-2) + curr, dfs(i-1))."
$ws.Range("A4").Value = "198. House Robber"
$ws.Range("B4").Value = "Medium"
$ws.Range("C4").Value = "Dynamic Programming"
$ws.Range("E4").Value = "https://leetcode.com/problems/house-robber/solutions/156523/from-good-to-great-how-to-approach-most-of-dp-problems/ "

# --- Hyperlinks for the two new Link cells ---
$ws.Hyperlinks.Add($ws.Range("E3"), "https://leetcode.com/problems/best-time-to-buy-and-sell-stock/solutions/1735493/java-c-best-ever-explanation-could-possible/")
$ws.Hyperlinks.Add($ws.Range("E4"), "https://leetcode.com/problems/house-robber/solutions/156523/from-good-to-great-how-to-approach-most-of-dp-problems/")

# Re-apply the hyperlink cell style (Hyperlinks.Add resets it to its own style index);
# restore it back to match the existing hyperlink-styled cell E2.
$ws.Range("E2").Copy()
$ws.Range("E3").PasteSpecial(-4122)
$ws.Range("E4").PasteSpecial(-4122)

# --- Resize table to include the two new rows ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E4"))

# --- Column C width grew to fit the new "Dynamic Programming" text ---
$ws.Columns.Item(3).AutoFit()

# --- Sheet view: scroll back to A1 and select D9 ---
$ws.Activate()
[void]$ws.Range("D9").Select()

Write-Host "done"
